# Automatische test-sync: 2025-06-23 18:09:50
# Adds the new "Open sollicitatie" mail-log entry (row 4) to the Logs
# sheet, the matching category tally row on the Dashboard sheet, and
# widens the conditional formatting + chart series ranges so they keep
# covering the newly added row.

$wb = $excel.ActiveWorkbook
$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A4:G4 -----------------------------------------------------
$wsLogs.Range("A4").Value = "Open sollicitatie"
$wsLogs.Range("B4").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C4").Value = "Zijn er op dit moment openstaande functies bij jullie bedrijf?"
$wsLogs.Range("D4").Value = "Sollicitatie / Vacature"
# (E4 intentionally left blank - no "Antwoord" recorded for this row)
$wsLogs.Range("F4").Value = "2025-06-23 18:09:41"
$wsLogs.Range("G4").Value = "Nee"

# Widen the conditional-formatting coverage from row 2:3 to 2:4 for
# both the Categorie (D) and Beantwoord (G) columns.
$fcD = $wsLogs.Range("D2:D3").FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($wsLogs.Range("D2:D4"))

$fcG = $wsLogs.Range("G2:G3").FormatConditions.Item(1)
$fcG.ModifyAppliesToRange($wsLogs.Range("G2:G4"))

# --- Dashboard!A4:B4 --------------------------------------------------
$wsDash.Range("A4").Value = "Sollicitatie / Vacature"
$wsDash.Range("B4").Value = 1

# Extend the bar chart's category/value series to include the new row.
$chartObj = $wsDash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
